# Paediatric motor neuronopathies.xlsx
# - Refresh the panel-query timestamps in the "data" sheet's time_taken column
# - Add a new "metadata" sheet (placed after "data") describing the panel query

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1. Update the F column (time_taken) timestamps for rows 2..40 -----------
$newTimes = @{
    2 = "2021-10-05 14:22:06.799255"
    3 = "2021-10-05 14:22:06.799264"
    4 = "2021-10-05 14:22:06.799267"
    5 = "2021-10-05 14:22:06.799270"
    6 = "2021-10-05 14:22:06.799273"
    7 = "2021-10-05 14:22:06.799275"
    8 = "2021-10-05 14:22:06.799278"
    9 = "2021-10-05 14:22:06.799281"
    10 = "2021-10-05 14:22:06.799284"
    11 = "2021-10-05 14:22:06.799287"
    12 = "2021-10-05 14:22:06.799289"
    13 = "2021-10-05 14:22:06.799292"
    14 = "2021-10-05 14:22:06.799294"
    15 = "2021-10-05 14:22:06.799297"
    16 = "2021-10-05 14:22:06.799300"
    17 = "2021-10-05 14:22:06.799302"
    18 = "2021-10-05 14:22:06.799305"
    19 = "2021-10-05 14:22:06.799308"
    20 = "2021-10-05 14:22:06.799310"
    21 = "2021-10-05 14:22:06.799313"
    22 = "2021-10-05 14:22:06.799316"
    23 = "2021-10-05 14:22:06.799318"
    24 = "2021-10-05 14:22:06.799321"
    25 = "2021-10-05 14:22:06.799323"
    26 = "2021-10-05 14:22:06.799326"
    27 = "2021-10-05 14:22:06.799329"
    28 = "2021-10-05 14:22:06.799332"
    29 = "2021-10-05 14:22:06.799334"
    30 = "2021-10-05 14:22:06.799337"
    31 = "2021-10-05 14:22:06.799339"
    32 = "2021-10-05 14:22:06.799342"
    33 = "2021-10-05 14:22:06.799345"
    34 = "2021-10-05 14:22:06.799348"
    35 = "2021-10-05 14:22:06.799350"
    36 = "2021-10-05 14:22:06.799353"
    37 = "2021-10-05 14:22:06.799356"
    38 = "2021-10-05 14:22:06.799358"
    39 = "2021-10-05 14:22:06.799361"
    40 = "2021-10-05 14:22:06.799363"
}

foreach ($row in $newTimes.Keys) {
    $ws.Cells.Item($row, 6).Value = $newTimes[$row]
}

# --- 2. Add the "metadata" sheet right after "data" ---------------------------
$metaSheet = $wb.Worksheets.Add($null, $ws)
$metaSheet.Name = "metadata"

# Copy the header formatting (bold + border + centered) from the data sheet's
# header row, and the index-column formatting from data!A2, onto the new sheet.
$ws.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# Header row
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Paediatric motor neuronopathies"
$metaSheet.Cells.Item(2, 3).Value = 79

# data_version ("1.69") must stay textual rather than be coerced to a number.
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.69"
$metaSheet.Cells.Item(2, 4).Style = "Normal"

$metaSheet.Cells.Item(2, 5).Value = "2021-07-09T13:03:54.797706Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:22:06.795543"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/79/?format=json"

# Keep "data" as the active sheet/tab (the workbook's bookViews are untouched
# by this change, so the active tab should stay where it was).
[void]$ws.Activate()
[void]$ws.Range("A1").Select()
